$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 307, shifting existing rows 307-410 down to 308-411.
$ws.Rows(307).Insert()

# Populate the newly inserted row 307 with the new record's data.
$ws.Range("A307").Value = 4
$ws.Range("B307").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C307").Value = "Los Lagos"
$ws.Range("D307").Value = 44988
$ws.Range("E307").Value = 10
$ws.Range("F307").Value = 100112040
$ws.Range("G307").Value = "Cilantro"
$ws.Range("H307").Value = "Sin especificar"
$ws.Range("I307").Value = "Primera"
$ws.Range("J307").Value = 90
$ws.Range("K307").Value = 9000
$ws.Range("L307").Value = 9000
$ws.Range("M307").Value = 9000
$ws.Range("N307").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O307").Value = "Región de La Araucanía"
$ws.Range("P307").Value = 4500
$ws.Range("Q307").Value = 2
$ws.Range("R307").Value = "Hortaliza"
